$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.015358805656433
$ws.Range("B1").Value = 1.193194150924683
$ws.Range("C1").Value = 5.550911903381348
$ws.Range("D1").Value = 1.638734340667725
$ws.Range("E1").Value = 1.001130342483521
